$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Empleados")

# The "Genero" column (B) contains values like " Masculino  ", " Femenino   ",
# " No Binario " etc. with stray leading/trailing spaces. Trim them so the
# column only contains clean "Masculino" / "Femenino" / "No Binario" values.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 50 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $trimmed = $val.ToString().Trim()
        if ($trimmed -ne $val) {
            $cell.Value2 = $trimmed
        }
    }
}

# Move the active selection to B2 (as seen after the edit).
$ws.Activate()
$ws.Range("B2").Select()
